# "minor grammatical updates to confusion matrix"
#
# Slide 3 contains the 3x3 confusion-matrix table ("Table 7"). A handful
# of its cells get bumped from 11pt to 14pt body text (the "No"/"Yes"
# row & column headers), and the four corner-result labels (TN / FP /
# FN / TP) get bolded for emphasis.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Find the confusion-matrix table shape on the slide.
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
    }
}

$tbl = $tableShape.Table

# Row 1: blank corner cell, "No" / "Yes" predicted-value headers.
$tbl.Cell(1, 2).Shape.TextFrame.TextRange.Font.Size = 14   # "No"
$tbl.Cell(1, 3).Shape.TextFrame.TextRange.Font.Size = 14   # "Yes"

# Row 2: "No" actual-value header, then TN / FP results.
$tbl.Cell(2, 1).Shape.TextFrame.TextRange.Font.Size = 14   # "No"
$tbl.Cell(2, 2).Shape.TextFrame.TextRange.Font.Bold = -1   # "TN"
$tbl.Cell(2, 3).Shape.TextFrame.TextRange.Font.Bold = -1   # "FP"

# Row 3: "Yes" actual-value header, then FN / TP results.
$tbl.Cell(3, 1).Shape.TextFrame.TextRange.Font.Size = 14   # "Yes"
$tbl.Cell(3, 2).Shape.TextFrame.TextRange.Font.Bold = -1   # "FN"
$tbl.Cell(3, 3).Shape.TextFrame.TextRange.Font.Bold = -1   # "TP"
